# Updates cryptos list data (prices and 1h volume change percentages)
# as scraped from coinranking.com, matching the GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the edited range to be treated as text so that values such as
# "504.00", "0.0000264" or "27.50" are not reinterpreted as numbers and
# lose their original textual formatting (trailing zeros, grouping dots, etc).
$editRange = $ws.Range("B2:E51")
$editRange.NumberFormat = "@"

$ws.Range("D2").Value = "65.002.29"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "3.147.57"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "579.58"
$ws.Range("E5").Value = "  +1.49%  "
$ws.Range("D6").Value = "148.87"
$ws.Range("E6").Value = "  -1.18%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.146.90"
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("D9").Value = "0.525"
$ws.Range("E9").Value = "  -0.77%  "
$ws.Range("E10").Value = "  -2.57%  "
$ws.Range("D11").Value = "6.14"
$ws.Range("E11").Value = "  -0.81%  "
$ws.Range("D12").Value = "0.499"
$ws.Range("E12").Value = "  -1.09%  "
$ws.Range("D13").Value = "0.0000264"
$ws.Range("E13").Value = "  +0.84%  "
$ws.Range("D14").Value = "37.14"
$ws.Range("E14").Value = "  -2.80%  "
$ws.Range("D15").Value = "3.660.64"
$ws.Range("E15").Value = "  -0.32%  "
$ws.Range("D16").Value = "64.911.15"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").Value = "7.14"
$ws.Range("E17").Value = "  -0.98%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.139.85"
$ws.Range("E18").Value = "  -0.47%  "
$ws.Range("D20").Value = "504.00"
$ws.Range("E20").Value = "  -2.38%  "
$ws.Range("D21").Value = "15.07"
$ws.Range("E21").Value = "  +0.95%  "
$ws.Range("D22").Value = "0.714"
$ws.Range("E22").Value = "  -3.11%  "
$ws.Range("D23").Value = "15.16"
$ws.Range("E23").Value = "  -1.79%  "
$ws.Range("D24").Value = "7.73"
$ws.Range("E24").Value = "  -1.83%  "
$ws.Range("D25").Value = "84.22"
$ws.Range("E25").Value = "  -0.97%  "
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("D27").Value = "9.12"
$ws.Range("E27").Value = "  +2.09%  "
$ws.Range("E28").Value = "  -0.55%  "
$ws.Range("E29").Value = "  -0.84%  "
$ws.Range("D30").Value = "2.79"
$ws.Range("E30").Value = "  +2.57%  "
$ws.Range("D31").Value = "27.50"
$ws.Range("E31").Value = "  -1.37%  "
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "6.43"
$ws.Range("E33").Value = "  +3.24%  "
$ws.Range("B34").Value = "Mantle"
$ws.Range("C34").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D34").Value = "1.19"
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("D35").Value = "6.48"
$ws.Range("E35").Value = "  -2.46%  "
$ws.Range("D36").Value = "54.94"
$ws.Range("D37").Value = "0.0887"
$ws.Range("E37").Value = "  +1.99%  "
$ws.Range("D38").Value = "475.54"
$ws.Range("E38").Value = "  -2.06%  "
$ws.Range("D39").Value = "0.0413"
$ws.Range("E39").Value = "  -2.69%  "
$ws.Range("D40").Value = "2.93"
$ws.Range("E40").Value = "  -4.15%  "
$ws.Range("D41").Value = "8.75"
$ws.Range("E41").Value = "  +0.87%  "
$ws.Range("D42").Value = "2.999.10"
$ws.Range("E42").Value = "  -3.91%  "
$ws.Range("E43").Value = "  -1.87%  "
$ws.Range("D44").Value = "0.282"
$ws.Range("E44").Value = "  -3.64%  "
$ws.Range("E45").Value = "  -2.82%  "
$ws.Range("D46").Value = "28.31"
$ws.Range("E46").Value = "  -4.66%  "
$ws.Range("D47").Value = "0.0₃0595"
$ws.Range("E47").Value = "  +2.89%  "
$ws.Range("E49").Value = "  -1.75%  "
$ws.Range("D50").Value = "2.25"
$ws.Range("E50").Value = "  -3.09%  "
$ws.Range("D51").Value = "2.48"
$ws.Range("E51").Value = "  +14.45%  "

# Restore default (unstyled) formatting now that the text values are set,
# so the cells keep matching the workbook's original (style-less) appearance.
$editRange.ClearFormats()
